# Daily refresh of the "剩余" (days remaining) / "开始时间" (cycle start date)
# tracker: advance the reference date to 2026-02-09 and recompute every data
# row (2..last) from its 总天 (D, total days) and 开始时间 (F, start date,
# stored as a YYYYMMDD integer):
#   end date   = F + D - 1 days
#   remaining  = end date - today + 1
#   if remaining <= 0  -> cycle finished: reset E to D (full cycle) and
#                         F to today (new cycle starts today)
#   else               -> E becomes the recomputed remaining, F unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$todaySerial = 20260209
$today = [datetime]::ParseExact([string]$todaySerial, "yyyyMMdd", $null)
$todayOA = $today.ToOADate()

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    if ($d -eq $null -or $f -eq $null) {
        continue
    }

    # Some rows carry a malformed start date (not a real calendar date) -
    # leave those rows untouched, same as the source edit did.
    $badDate = $false
    try {
        $fdt = [datetime]::ParseExact([string]$f, "yyyyMMdd", $null)
    } catch {
        $badDate = $true
    }

    if ($badDate) {
        continue
    }

    $endOA = $fdt.ToOADate() + $d - 1
    $remaining = $endOA - $todayOA + 1

    if ($remaining -le 0) {
        $newE = $d
        $newF = $todaySerial
    } else {
        $newE = $remaining
        $newF = $f
    }

    $ws.Cells.Item($r, 5).Value = $newE
    $ws.Cells.Item($r, 6).Value = $newF
}
